$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H32").Value = 5221.727
$ws.Range("J32").Value = 5348.4287
$ws.Range("L32").Value = 5348.4287
$ws.Range("N32").Value = -6000.4287
$ws.Range("H48").Value = 2789.4285
$ws.Range("I48").Value = 3025.2
$ws.Range("J48").Value = 2200
$ws.Range("K48").Value = 9075.599999999999
$ws.Range("L48").Value = 6600
$ws.Range("M48").Value = -8783.599999999999
$ws.Range("N48").Value = -7184
$ws.Range("H56").Value = 2789.4285
$ws.Range("I56").Value = 3025.2
$ws.Range("J56").Value = 2200
$ws.Range("K56").Value = 9075.599999999999
$ws.Range("L56").Value = 6600
$ws.Range("M56").Value = -8541.599999999999
$ws.Range("N56").Value = -7668
$ws.Range("H106").Value = 2413.3333
$ws.Range("I106").Value = 1959.1428
$ws.Range("K106").Value = 1959.1428
$ws.Range("M106").Value = -1328.1428
$ws.Range("H111").Value = 1747.5
$ws.Range("I111").Value = 536.8889
$ws.Range("J111").Value = 2958.111
$ws.Range("K111").Value = 1610.6667
$ws.Range("L111").Value = 8874.332999999999
$ws.Range("M111").Value = 1456.3333
$ws.Range("N111").Value = -15008.333
$ws.Range("H129").Value = 1856.6897
$ws.Range("I129").Value = 1944.5834
$ws.Range("J129").Value = 1434.8
$ws.Range("K129").Value = 5833.7502
$ws.Range("L129").Value = 4304.4
$ws.Range("M129").Value = -833.7502000000004
$ws.Range("N129").Value = -14304.4

$ws = $wb.Worksheets("ARM")
$ws.Range("H2").Value = 592.0345
$ws.Range("I2").Value = 530.55
$ws.Range("J2").Value = 728.6667
$ws.Range("K2").Value = 530.55
$ws.Range("L2").Value = 728.6667
$ws.Range("M2").Value = -417.55
$ws.Range("N2").Value = -954.6667
$ws.Range("H19").Value = 30999.666
$ws.Range("I19").Value = 30999.666
$ws.Range("K19").Value = 30999.666
$ws.Range("M19").Value = -30770.666
$ws.Range("H32").Value = 10593.392
$ws.Range("I32").Value = 5801.5
$ws.Range("K32").Value = 5801.5
$ws.Range("M32").Value = -5514.5
$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -50976
$ws.Range("H55").Value = 14599.6
$ws.Range("I55").Value = 11000
$ws.Range("K55").Value = 11000
$ws.Range("M55").Value = -10685
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = 0
$ws.Range("H110").Value = 7184.7144
$ws.Range("I110").Value = 7230.4736
$ws.Range("K110").Value = 7230.4736
$ws.Range("M110").Value = -5185.4736
$ws.Range("H116").Value = 592.0345
$ws.Range("I116").Value = 530.55
$ws.Range("J116").Value = 728.6667
$ws.Range("K116").Value = 530.55
$ws.Range("L116").Value = 728.6667
$ws.Range("M116").Value = 1763.45
$ws.Range("N116").Value = -5316.6667
$ws.Range("H132").Value = 1798.907
$ws.Range("I132").Value = 1545.317
$ws.Range("K132").Value = 4635.951
$ws.Range("M132").Value = -2105.951

$ws = $wb.Worksheets("BSM")
$ws.Range("H3").Value = 592.0345
$ws.Range("I3").Value = 530.55
$ws.Range("J3").Value = 728.6667
$ws.Range("K3").Value = 530.55
$ws.Range("L3").Value = 728.6667
$ws.Range("M3").Value = -416.55
$ws.Range("N3").Value = -956.6667
$ws.Range("H20").Value = 2968.6606
$ws.Range("I20").Value = 2535.5757
$ws.Range("J20").Value = 3590.0435
$ws.Range("K20").Value = 2535.5757
$ws.Range("L20").Value = 3590.0435
$ws.Range("M20").Value = -2288.5757
$ws.Range("N20").Value = -4084.0435
$ws.Range("H94").Value = 1409.0741
$ws.Range("I94").Value = 1245.3334
$ws.Range("K94").Value = 1245.3334
$ws.Range("M94").Value = -794.3334
$ws.Range("H105").Value = 2590
$ws.Range("I105").Value = 1860.1666
$ws.Range("J105").Value = 6969
$ws.Range("K105").Value = 1860.1666
$ws.Range("L105").Value = 6969
$ws.Range("M105").Value = -113.1666
$ws.Range("N105").Value = -10463
$ws.Range("H134").Value = 1882.5942
$ws.Range("I134").Value = 1896.2538
$ws.Range("K134").Value = 5688.761399999999
$ws.Range("M134").Value = -3153.761399999999

$ws = $wb.Worksheets("CRP")
$ws.Range("H99").Value = 3643.75
$ws.Range("I99").Value = 3419.2
$ws.Range("K99").Value = 3419.2
$ws.Range("M99").Value = -1921.2
$ws.Range("H107").Value = 969.17645
$ws.Range("I107").Value = 599
$ws.Range("K107").Value = 599
$ws.Range("M107").Value = 1321
$ws.Range("H126").Value = 3643.75
$ws.Range("I126").Value = 3419.2
$ws.Range("K126").Value = 10257.6
$ws.Range("M126").Value = -7787.599999999999
$ws.Range("H127").Value = 100769.336
$ws.Range("J127").Value = 100769.336
$ws.Range("L127").Value = 100769.336
$ws.Range("N127").Value = -110689.336
$ws.Range("H132").Value = 4552

$ws = $wb.Worksheets("CUL")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864

$ws = $wb.Worksheets("GSM")
$ws.Range("H113").Value = 2959.9285
$ws.Range("I113").Value = 2959.9285
$ws.Range("K113").Value = 2959.9285
$ws.Range("M113").Value = -789.9285

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 4035.3845
$ws.Range("I7").Value = 3878.55
$ws.Range("K7").Value = 3878.55
$ws.Range("M7").Value = -3766.55
$ws.Range("H22").Value = 2326.8572
$ws.Range("I22").Value = 2192
$ws.Range("J22").Value = 2428
$ws.Range("K22").Value = 2192
$ws.Range("L22").Value = 2428
$ws.Range("M22").Value = -1897
$ws.Range("N22").Value = -3018
$ws.Range("H27").Value = 2326.8572
$ws.Range("I27").Value = 2192
$ws.Range("J27").Value = 2428
$ws.Range("K27").Value = 2192
$ws.Range("L27").Value = 2428
$ws.Range("M27").Value = -2085
$ws.Range("N27").Value = -2642
$ws.Range("H46").Value = 1571
$ws.Range("I46").Value = 1100
$ws.Range("K46").Value = 1100
$ws.Range("M46").Value = -912
$ws.Range("H55").Value = 105.09524
$ws.Range("I55").Value = 111.68421
$ws.Range("J55").Value = 42.5
$ws.Range("K55").Value = 111.68421
$ws.Range("L55").Value = 42.5
$ws.Range("M55").Value = 61.31579000000001
$ws.Range("N55").Value = -388.5
$ws.Range("H126").Value = 4035.3845
$ws.Range("I126").Value = 3878.55
$ws.Range("K126").Value = 11635.65
$ws.Range("M126").Value = -9165.650000000001
$ws.Range("H132").Value = 4384.5
$ws.Range("I132").Value = 3924.2144
$ws.Range("K132").Value = 11772.6432
$ws.Range("M132").Value = -9242.643199999999
$ws.Range("H136").Value = 3789.7097
$ws.Range("I136").Value = 3419.4
$ws.Range("K136").Value = 10258.2
$ws.Range("M136").Value = -7708.200000000001

$ws = $wb.Worksheets("WVR")
$ws.Range("H17").Value = 54999.5
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H43").Value = 21756.75
$ws.Range("I43").Value = 12027
$ws.Range("J43").Value = 25000
$ws.Range("K43").Value = 12027
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = -11878
$ws.Range("N43").Value = -25298
$ws.Range("H81").Value = 9969.529
$ws.Range("I81").Value = 51598
$ws.Range("J81").Value = 4419.067
$ws.Range("K81").Value = 103196
$ws.Range("L81").Value = 8838.134
$ws.Range("M81").Value = -102135
$ws.Range("N81").Value = -10960.134
$ws.Range("H84").Value = 9969.529
$ws.Range("I84").Value = 51598
$ws.Range("J84").Value = 4419.067
$ws.Range("K84").Value = 515980
$ws.Range("L84").Value = 44190.67
$ws.Range("M84").Value = -510676
$ws.Range("N84").Value = -54798.67
$ws.Range("H126").Value = 3879.2307
$ws.Range("I126").Value = 3418.7144
$ws.Range("J126").Value = 4416.5
$ws.Range("K126").Value = 10256.1432
$ws.Range("L126").Value = 13249.5
$ws.Range("M126").Value = -7786.143199999999
$ws.Range("N126").Value = -18189.5
$ws.Range("H132").Value = 2064.0967
$ws.Range("I132").Value = 2099.261
$ws.Range("K132").Value = 6297.782999999999
$ws.Range("M132").Value = -3767.782999999999
$ws.Range("H139").Value = 119959.75
$ws.Range("J139").Value = 126613
$ws.Range("L139").Value = 126613
$ws.Range("N139").Value = -136893
